# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Sat Aug 26 15:53:41 UTC 2023 with GitHub Actions")
#
# Numeric-looking price strings (column D) are prefixed with a leading
# apostrophe so Excel keeps storing them as text, matching the original
# inlineStr cell type instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.147.77'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.656.29'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''218.09'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').Value = '''0.5298'
$ws.Range('E6').Value = '  +1.78%  '
$ws.Range('D7').Value = '''1.004'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '''0.2614'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''0.06333'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').Value = '''20.46'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').Value = '''0.07809'
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').Value = '''4.519'
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').Value = '1.632.96'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').Value = '1.883.91'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').Value = '''0.5496'
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').Value = '0.0₅8215'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').Value = '''65.42'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '26.148.80'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').Value = '''1.004'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '''4.605'
$ws.Range('E20').Value = '  +0.81%  '
$ws.Range('D21').Value = '''191.49'
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('D23').Value = '''6.021'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = '''144.72'
$ws.Range('E25').Value = '  +4.43%  '
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('D27').Value = '''7.220'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('D28').Value = '''16.00'
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('E29').Value = '  +4.29%  '
$ws.Range('D30').Value = '''0.05747'
$ws.Range('E30').Value = '  -3.66%  '
$ws.Range('D31').Value = '''1.275'
$ws.Range('E31').Value = '  +0.26%  '
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('D33').Value = '''3.274'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').Value = '''1.602'
$ws.Range('E34').Value = '  +3.43%  '
$ws.Range('D35').Value = '''2.805'
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('D36').Value = '''0.9524'
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('D37').Value = '''2.419'
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').Value = '''0.5747'
$ws.Range('E38').Value = '  +0.96%  '
$ws.Range('D39').Value = '''0.01611'
$ws.Range('E39').Value = '  +0.81%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''0.8521'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''5.802'
$ws.Range('E41').Value = '  -1.37%  '
$ws.Range('D42').Value = '''104.50'
$ws.Range('E42').Value = '  +3.64%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '''1.004'
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.042.28'
$ws.Range('E44').Value = '  +3.65%  '
$ws.Range('D45').Value = '1.797.78'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').Value = '''56.93'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''7.881'
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.05156'
$ws.Range('E50').Value = '  +0.07%  '
$ws.Range('D51').Value = '''1.445'
$ws.Range('E51').Value = '  -2.51%  '
